$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 115
$ws.Range("I9").Value = 146.5
$ws.Range("J9").Value = 52
$ws.Range("K9").Value = 146.5
$ws.Range("L9").Value = 52
$ws.Range("M9").Value = 22.5
$ws.Range("N9").Value = -390

$ws.Range("H21").Value = 34000
$ws.Range("J21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("N21").Value = -30936

$ws.Range("H23").Value = 34000
$ws.Range("J23").Value = 30000
$ws.Range("L23").Value = 30000
$ws.Range("N23").Value = -30468

$ws.Range("H33").Value = 265.4
$ws.Range("I33").Value = 283.77777
$ws.Range("K33").Value = 283.77777
$ws.Range("M33").Value = -54.77776999999998

$ws.Range("H38").Value = 386.81818
$ws.Range("J38").Value = 3500
$ws.Range("L38").Value = 10500
$ws.Range("N38").Value = -11244

$ws.Range("H113").Value = 2976.25
$ws.Range("I113").Value = 2968.3333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2968.3333
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 285.6667000000002
$ws.Range("N113").Value = -9508

$ws.Range("H138").Value = 4440.98
$ws.Range("I138").Value = 2520.4546
$ws.Range("J138").Value = 6788.289
$ws.Range("K138").Value = 7561.3638
$ws.Range("L138").Value = 20364.867
$ws.Range("M138").Value = -2421.3638
$ws.Range("N138").Value = -30644.867

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 143916.14
$ws.Range("I2").Value = 1100
$ws.Range("J2").Value = 334337.66
$ws.Range("K2").Value = 1100
$ws.Range("L2").Value = 334337.66
$ws.Range("M2").Value = -987
$ws.Range("N2").Value = -334563.66

$ws.Range("H32").Value = 33635.344
$ws.Range("I32").Value = 28053.24
$ws.Range("K32").Value = 28053.24
$ws.Range("M32").Value = -27766.24

$ws.Range("H45").Value = 50606
$ws.Range("I45").Value = 1212
$ws.Range("J45").Value = 100000
$ws.Range("K45").Value = 1212
$ws.Range("L45").Value = 100000
$ws.Range("M45").Value = -835
$ws.Range("N45").Value = -100754

$ws.Range("H76").Value = 124666.336
$ws.Range("J76").Value = 124666.336
$ws.Range("L76").Value = 124666.336
$ws.Range("N76").Value = -125342.336

$ws.Range("H79").Value = 124666.336
$ws.Range("J79").Value = 124666.336
$ws.Range("L79").Value = 124666.336
$ws.Range("N79").Value = -127006.336

$ws.Range("H105").Value = 275185
$ws.Range("J105").Value = 275185
$ws.Range("L105").Value = 275185
$ws.Range("N105").Value = -282173

$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = $null

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null

$ws.Range("H116").Value = 143916.14
$ws.Range("I116").Value = 1100
$ws.Range("J116").Value = 334337.66
$ws.Range("K116").Value = 1100
$ws.Range("L116").Value = 334337.66
$ws.Range("M116").Value = 1194
$ws.Range("N116").Value = -338925.66

$ws.Range("H122").Value = 2455.6365
$ws.Range("I122").Value = 2334.6667
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7004.000100000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4554.000100000001
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 143916.14
$ws.Range("I3").Value = 1100
$ws.Range("J3").Value = 334337.66
$ws.Range("K3").Value = 1100
$ws.Range("L3").Value = 334337.66
$ws.Range("M3").Value = -986
$ws.Range("N3").Value = -334565.66

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = $null

$ws.Range("H64").Value = 235.78572
$ws.Range("I64").Value = 293.77777
$ws.Range("K64").Value = 293.77777
$ws.Range("M64").Value = -68.77776999999998

$ws.Range("H67").Value = 235.78572
$ws.Range("I67").Value = 293.77777
$ws.Range("K67").Value = 293.77777
$ws.Range("M67").Value = 486.22223

$ws.Range("H76").Value = 30519.8
$ws.Range("J76").Value = 34328.5
$ws.Range("L76").Value = 34328.5
$ws.Range("N76").Value = -34958.5

$ws.Range("H79").Value = 30519.8
$ws.Range("J79").Value = 34328.5
$ws.Range("L79").Value = 34328.5
$ws.Range("N79").Value = -36512.5

$ws.Range("H107").Value = 2143.3103
$ws.Range("I107").Value = 2028
$ws.Range("K107").Value = 2028
$ws.Range("M107").Value = -108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2475
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 2950
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2950
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -3524

$ws.Range("H58").Value = 1433.2188
$ws.Range("I58").Value = 1394.9667
$ws.Range("K58").Value = 1394.9667
$ws.Range("M58").Value = -1191.9667

$ws.Range("H113").Value = 2475
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 2950
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 2950
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -7290

$ws.Range("H132").Value = 1622.1143
$ws.Range("I132").Value = 1414.8438
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 4244.5314
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -1714.5314
$ws.Range("N132").Value = -16559

$ws.Range("H136").Value = 1433.2188
$ws.Range("I136").Value = 1394.9667
$ws.Range("K136").Value = 4184.9001
$ws.Range("M136").Value = -1634.9001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1017082.4
$ws.Range("I12").Value = 28.333334
$ws.Range("J12").Value = 1486491.9
$ws.Range("K12").Value = 85.00000199999999
$ws.Range("L12").Value = 4459475.699999999
$ws.Range("M12").Value = 87.99999800000001
$ws.Range("N12").Value = -4459821.699999999

$ws.Range("H34").Value = 934.8889
$ws.Range("J34").Value = 1092.8636
$ws.Range("L34").Value = 3278.5908
$ws.Range("N34").Value = -3446.5908

$ws.Range("H39").Value = 2937.7
$ws.Range("J39").Value = 2937.7
$ws.Range("L39").Value = 8813.099999999999
$ws.Range("N39").Value = -9401.099999999999

$ws.Range("H55").Value = 3339.8
$ws.Range("J55").Value = 3339.8
$ws.Range("L55").Value = 10019.4
$ws.Range("N55").Value = -10373.4

$ws.Range("H131").Value = 867.17
$ws.Range("J131").Value = 883.69147
$ws.Range("L131").Value = 2651.07441
$ws.Range("N131").Value = -12731.07441

$ws.Range("H136").Value = 4655.4165
$ws.Range("I136").Value = 1345.7142
$ws.Range("J136").Value = 6018.2354
$ws.Range("K136").Value = 4037.1426
$ws.Range("L136").Value = 18054.7062
$ws.Range("M136").Value = 1062.8574
$ws.Range("N136").Value = -28254.7062

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3552.125
$ws.Range("J43").Value = 6000
$ws.Range("L43").Value = 6000
$ws.Range("N43").Value = -6302

$ws.Range("H49").Value = 15955
$ws.Range("J49").Value = 15955
$ws.Range("L49").Value = 15955
$ws.Range("N49").Value = -16323

$ws.Range("H108").Value = 40684
$ws.Range("J108").Value = 40684
$ws.Range("L108").Value = 40684
$ws.Range("N108").Value = -48364

$ws.Range("H113").Value = 1979.4286
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 2109.3333
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 2109.3333
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -6449.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 19272.727
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20348

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").Value = $null

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null

$ws.Range("H122").Value = 10353984
$ws.Range("I122").Value = 13975778
$ws.Range("J122").Value = 6000.4287
$ws.Range("K122").Value = 41927334
$ws.Range("L122").Value = 18001.2861
$ws.Range("M122").Value = -41924884
$ws.Range("N122").Value = -22901.2861

$ws.Range("H132").Value = 6500
$ws.Range("I132").Value = 6500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 19500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -16970
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 35716064
$ws.Range("I122").Value = 41668076
$ws.Range("K122").Value = 125004228
$ws.Range("M122").Value = -125001778

$ws.Range("H132").Value = 1860.079
$ws.Range("I132").Value = 1512.3871
$ws.Range("J132").Value = 3399.8572
$ws.Range("K132").Value = 4537.1613
$ws.Range("L132").Value = 10199.5716
$ws.Range("M132").Value = -2007.1613
$ws.Range("N132").Value = -15259.5716

$ws.Range("H133").Value = 55224.332
$ws.Range("J133").Value = 55224.332
$ws.Range("L133").Value = 55224.332
$ws.Range("N133").Value = -65344.332
